$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 40.75339133333333
$ws.Cells.Item(2, 8).Value = 122.260174
$ws.Cells.Item(2, 9).Value = 0.02126536631186857
$ws.Cells.Item(2, 10).Value = 0.02126536631186857
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 15.35884066666667
$ws.Cells.Item(2, 14).Value = 46.076522
$ws.Cells.Item(2, 15).Value = 0.1012042817263867
$ws.Cells.Item(2, 16).Value = 0.1012042817263867
$ws.Cells.Item(2, 17).Value = 625.9248441149808
$ws.Cells.Item(2, 18).Value = 5633.323597034828
$ws.Cells.Item(2, 19).Value = 0.002152146123241158
$ws.Cells.Item(2, 20).Value = 0.002152146123241159

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 40.75339133333333
$ws.Cells.Item(3, 8).Value = 122.260174
$ws.Cells.Item(3, 9).Value = 0.02126536631186857
$ws.Cells.Item(3, 10).Value = 0.02126536631186857
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 50.59256466666667
$ws.Cells.Item(3, 14).Value = 151.777694
$ws.Cells.Item(3, 15).Value = 0.3333704853712116
$ws.Cells.Item(3, 16).Value = 0.3333704853712116
$ws.Cells.Item(3, 17).Value = 2061.818586417639
$ws.Cells.Item(3, 18).Value = 18556.36727775876
$ws.Cells.Item(3, 19).Value = 0.007089245488984235
$ws.Cells.Item(3, 20).Value = 0.007089245488984235

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 40.75339133333333
$ws.Cells.Item(4, 8).Value = 122.260174
$ws.Cells.Item(4, 9).Value = 0.02126536631186857
$ws.Cells.Item(4, 10).Value = 0.02126536631186857
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 60.37715666666667
$ws.Cells.Item(4, 14).Value = 181.13147
$ws.Cells.Item(4, 15).Value = 0.397844271305776
$ws.Cells.Item(4, 16).Value = 0.397844271305776
$ws.Cells.Item(4, 17).Value = 2460.573893230643
$ws.Cells.Item(4, 18).Value = 22145.16503907578
$ws.Cells.Item(4, 19).Value = 0.008460304164395747
$ws.Cells.Item(4, 20).Value = 0.008460304164395747

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 40.75339133333333
$ws.Cells.Item(5, 8).Value = 122.260174
$ws.Cells.Item(5, 9).Value = 0.02126536631186857
$ws.Cells.Item(5, 10).Value = 0.02126536631186857
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 25.43221733333333
$ws.Cells.Item(5, 14).Value = 76.29665199999999
$ws.Cells.Item(5, 15).Value = 0.1675809615966257
$ws.Cells.Item(5, 16).Value = 0.1675809615966258
$ws.Cells.Item(5, 17).Value = 1036.449105459716
$ws.Cells.Item(5, 18).Value = 9328.041949137449
$ws.Cells.Item(5, 19).Value = 0.003563670535247425
$ws.Cells.Item(5, 20).Value = 0.003563670535247426

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1689.289306666667
$ws.Cells.Item(6, 8).Value = 5067.86792
$ws.Cells.Item(6, 9).Value = 0.8814813868902838
$ws.Cells.Item(6, 10).Value = 0.8814813868902838
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 15.35884066666667
$ws.Cells.Item(6, 14).Value = 46.076522
$ws.Cells.Item(6, 15).Value = 0.1012042817263867
$ws.Cells.Item(6, 16).Value = 0.1012042817263867
$ws.Cells.Item(6, 17).Value = 25945.52530099713
$ws.Cells.Item(6, 18).Value = 233509.7277089742
$ws.Cells.Item(6, 19).Value = 0.08920969061541031
$ws.Cells.Item(6, 20).Value = 0.08920969061541031

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1689.289306666667
$ws.Cells.Item(7, 8).Value = 5067.86792
$ws.Cells.Item(7, 9).Value = 0.8814813868902838
$ws.Cells.Item(7, 10).Value = 0.8814813868902838
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 50.59256466666667
$ws.Cells.Item(7, 14).Value = 151.777694
$ws.Cells.Item(7, 15).Value = 0.3333704853712116
$ws.Cells.Item(7, 16).Value = 0.3333704853712116
$ws.Cells.Item(7, 17).Value = 85465.47848824183
$ws.Cells.Item(7, 18).Value = 769189.3063941764
$ws.Cells.Item(7, 19).Value = 0.2938598777933026
$ws.Cells.Item(7, 20).Value = 0.2938598777933026

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1689.289306666667
$ws.Cells.Item(8, 8).Value = 5067.86792
$ws.Cells.Item(8, 9).Value = 0.8814813868902838
$ws.Cells.Item(8, 10).Value = 0.8814813868902838
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 60.37715666666667
$ws.Cells.Item(8, 14).Value = 181.13147
$ws.Cells.Item(8, 15).Value = 0.397844271305776
$ws.Cells.Item(8, 16).Value = 0.397844271305776
$ws.Cells.Item(8, 17).Value = 101994.4851239381
$ws.Cells.Item(8, 18).Value = 917950.3661154424
$ws.Cells.Item(8, 19).Value = 0.3506923200369698
$ws.Cells.Item(8, 20).Value = 0.3506923200369698

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1689.289306666667
$ws.Cells.Item(9, 8).Value = 5067.86792
$ws.Cells.Item(9, 9).Value = 0.8814813868902838
$ws.Cells.Item(9, 10).Value = 0.8814813868902838
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 25.43221733333333
$ws.Cells.Item(9, 14).Value = 76.29665199999999
$ws.Cells.Item(9, 15).Value = 0.1675809615966257
$ws.Cells.Item(9, 16).Value = 0.1675809615966258
$ws.Cells.Item(9, 17).Value = 42962.37278602264
$ws.Cells.Item(9, 18).Value = 386661.3550742038
$ws.Cells.Item(9, 19).Value = 0.1477194984446011
$ws.Cells.Item(9, 20).Value = 0.1477194984446011

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 100.9654023333333
$ws.Cells.Item(10, 8).Value = 302.896207
$ws.Cells.Item(10, 9).Value = 0.05268435816499466
$ws.Cells.Item(10, 10).Value = 0.05268435816499466
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 15.35884066666667
$ws.Cells.Item(10, 14).Value = 46.076522
$ws.Cells.Item(10, 15).Value = 0.1012042817263867
$ws.Cells.Item(10, 16).Value = 0.1012042817263867
$ws.Cells.Item(10, 17).Value = 1550.711527283562
$ws.Cells.Item(10, 18).Value = 13956.40374555205
$ws.Cells.Item(10, 19).Value = 0.005331882626303978
$ws.Cells.Item(10, 20).Value = 0.005331882626303979

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 100.9654023333333
$ws.Cells.Item(11, 8).Value = 302.896207
$ws.Cells.Item(11, 9).Value = 0.05268435816499466
$ws.Cells.Item(11, 10).Value = 0.05268435816499466
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 50.59256466666667
$ws.Cells.Item(11, 14).Value = 151.777694
$ws.Cells.Item(11, 15).Value = 0.3333704853712116
$ws.Cells.Item(11, 16).Value = 0.3333704853712116
$ws.Cells.Item(11, 17).Value = 5108.098646645184
$ws.Cells.Item(11, 18).Value = 45972.88781980666
$ws.Cells.Item(11, 19).Value = 0.01756341005293502
$ws.Cells.Item(11, 20).Value = 0.01756341005293502

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 100.9654023333333
$ws.Cells.Item(12, 8).Value = 302.896207
$ws.Cells.Item(12, 9).Value = 0.05268435816499466
$ws.Cells.Item(12, 10).Value = 0.05268435816499466
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 60.37715666666667
$ws.Cells.Item(12, 14).Value = 181.13147
$ws.Cells.Item(12, 15).Value = 0.397844271305776
$ws.Cells.Item(12, 16).Value = 0.397844271305776
$ws.Cells.Item(12, 17).Value = 6096.003914592699
$ws.Cells.Item(12, 18).Value = 54864.0352313343
$ws.Cells.Item(12, 19).Value = 0.02096017008336481
$ws.Cells.Item(12, 20).Value = 0.02096017008336481

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 100.9654023333333
$ws.Cells.Item(13, 8).Value = 302.896207
$ws.Cells.Item(13, 9).Value = 0.05268435816499466
$ws.Cells.Item(13, 10).Value = 0.05268435816499466
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 25.43221733333333
$ws.Cells.Item(13, 14).Value = 76.29665199999999
$ws.Cells.Item(13, 15).Value = 0.1675809615966257
$ws.Cells.Item(13, 16).Value = 0.1675809615966258
$ws.Cells.Item(13, 17).Value = 2567.774055288774
$ws.Cells.Item(13, 18).Value = 23109.96649759896
$ws.Cells.Item(13, 19).Value = 0.008828895402390846
$ws.Cells.Item(13, 20).Value = 0.008828895402390848

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 85.41274733333334
$ws.Cells.Item(14, 8).Value = 256.238242
$ws.Cells.Item(14, 9).Value = 0.04456888863285297
$ws.Cells.Item(14, 10).Value = 0.04456888863285297
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 15.35884066666667
$ws.Cells.Item(14, 14).Value = 46.076522
$ws.Cells.Item(14, 15).Value = 0.1012042817263867
$ws.Cells.Item(14, 16).Value = 0.1012042817263867
$ws.Cells.Item(14, 17).Value = 1311.840777194925
$ws.Cells.Item(14, 18).Value = 11806.56699475432
$ws.Cells.Item(14, 19).Value = 0.004510562361431203
$ws.Cells.Item(14, 20).Value = 0.004510562361431204

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 85.41274733333334
$ws.Cells.Item(15, 8).Value = 256.238242
$ws.Cells.Item(15, 9).Value = 0.04456888863285297
$ws.Cells.Item(15, 10).Value = 0.04456888863285297
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 50.59256466666667
$ws.Cells.Item(15, 14).Value = 151.777694
$ws.Cells.Item(15, 15).Value = 0.3333704853712116
$ws.Cells.Item(15, 16).Value = 0.3333704853712116
$ws.Cells.Item(15, 17).Value = 4321.249942819328
$ws.Cells.Item(15, 18).Value = 38891.24948537395
$ws.Cells.Item(15, 19).Value = 0.01485795203598967
$ws.Cells.Item(15, 20).Value = 0.01485795203598967

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 85.41274733333334
$ws.Cells.Item(16, 8).Value = 256.238242
$ws.Cells.Item(16, 9).Value = 0.04456888863285297
$ws.Cells.Item(16, 10).Value = 0.04456888863285297
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 60.37715666666667
$ws.Cells.Item(16, 14).Value = 181.13147
$ws.Cells.Item(16, 15).Value = 0.397844271305776
$ws.Cells.Item(16, 16).Value = 0.397844271305776
$ws.Cells.Item(16, 17).Value = 5156.978827075083
$ws.Cells.Item(16, 18).Value = 46412.80944367575
$ws.Cells.Item(16, 19).Value = 0.01773147702104567
$ws.Cells.Item(16, 20).Value = 0.01773147702104567

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 85.41274733333334
$ws.Cells.Item(17, 8).Value = 256.238242
$ws.Cells.Item(17, 9).Value = 0.04456888863285297
$ws.Cells.Item(17, 10).Value = 0.04456888863285297
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 25.43221733333333
$ws.Cells.Item(17, 14).Value = 76.29665199999999
$ws.Cells.Item(17, 15).Value = 0.1675809615966257
$ws.Cells.Item(17, 16).Value = 0.1675809615966258
$ws.Cells.Item(17, 17).Value = 2172.235553218421
$ws.Cells.Item(17, 18).Value = 19550.11997896578
$ws.Cells.Item(17, 19).Value = 0.007468897214386423
$ws.Cells.Item(17, 20).Value = 0.007468897214386424
